# Generate Report for Handback
# Reproduces: marking handoff rows as handed-back, recording handback file
# names/timestamps, and widening columns to fit the new content.

$wb = $excel.ActiveWorkbook

$mdUrl26 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1f3e515cd41cd68962877e56ca22dbb4e75ed5e/e2e/26c8ffe9-8196-4ba6-ab2c-3768decddbe4.md"
$mdUrl48 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1f3e515cd41cd68962877e56ca22dbb4e75ed5e/e2e/484e50fd-89bc-4ee3-a875-834b42a23e30.md"

# ---------------------------------------------------------------------
# Overview sheet: the "handoff status" columns mirror the per-language
# sheets' Status column text, so they get the same new status string.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666666666667

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 -> 26c8ffe9...
$wsZh.Range("I2").Value = "26c8ffe9-8196-4ba6-ab2c-3768decddbe4.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl26, "", "", "26c8ffe9-8196-4ba6-ab2c-3768decddbe4.md")
$wsZh.Range("J2").Value = "26c8ffe9-8196-4ba6-ab2c-3768decddbe4.9cf9ce5121b347886d8a0a7be649744db25b3a56.zh-cn.xlf"

# Row 3 -> 484e50fd...
$wsZh.Range("I3").Value = "484e50fd-89bc-4ee3-a875-834b42a23e30.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl48, "", "", "484e50fd-89bc-4ee3-a875-834b42a23e30.md")
$wsZh.Range("J3").Value = "484e50fd-89bc-4ee3-a875-834b42a23e30.0d533791aeac50fb6e6f0fc632145821d8ee54e9.zh-cn.xlf"

# Latest Handback DateTime already holds the placeholder "0001-01-01 00:00:00";
# stamp it with the actual handback time now that the target file is in sync.
$wsZh.Range("K2").Value = "2016-09-03 12:32:41"
$wsZh.Range("K3").Value = "2016-09-03 12:32:41"

$wsZh.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsZh.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsZh.Columns.Item(10).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2 -> 26c8ffe9...
$wsDe.Range("I2").Value = "26c8ffe9-8196-4ba6-ab2c-3768decddbe4.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl26, "", "", "26c8ffe9-8196-4ba6-ab2c-3768decddbe4.md")
$wsDe.Range("J2").Value = "26c8ffe9-8196-4ba6-ab2c-3768decddbe4.9cf9ce5121b347886d8a0a7be649744db25b3a56.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 12:32:47"

# Row 3 -> 484e50fd...
$wsDe.Range("I3").Value = "484e50fd-89bc-4ee3-a875-834b42a23e30.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl48, "", "", "484e50fd-89bc-4ee3-a875-834b42a23e30.md")
$wsDe.Range("J3").Value = "484e50fd-89bc-4ee3-a875-834b42a23e30.0d533791aeac50fb6e6f0fc632145821d8ee54e9.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 12:32:47"

$wsDe.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsDe.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsDe.Columns.Item(10).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# (set on both language sheets' Status column; Overview!E/F were updated
# above with the matching text.)
# ---------------------------------------------------------------------
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

Write-Host "Handback report generated"
